{"js": "// Replace the date line and the 25 \"two-digit \u00d7 two-digit\" answer cells\n// with the updated values from the new day's worksheet. Every old value\n// in this document is unique, so a simple exact, case-sensitive\n// search-and-replace (first match) for each pair is unambiguous.\nconst replacements = [\n  [\"2024-12-11 Wednesday\", \"2024-12-12 Thursday\"],\n  [\"44\u00d728=1232\", \"74\u00d740=2960\"],\n  [\"78\u00d728=2184\", \"23\u00d724=552\"],\n  [\"22\u00d795=2090\", \"87\u00d718=1566\"],\n  [\"41\u00d734=1394\", \"73\u00d775=5475\"],\n  [\"19\u00d723=437\", \"83\u00d778=6474\"],\n  [\"66\u00d723=1518\", \"94\u00d717=1598\"],\n  [\"39\u00d726=1014\", \"40\u00d779=3160\"],\n  [\"33\u00d775=2475\", \"13\u00d791=1183\"],\n  [\"98\u00d756=5488\", \"86\u00d760=5160\"],\n  [\"79\u00d721=1659\", \"39\u00d728=1092\"],\n  [\"48\u00d769=3312\", \"34\u00d786=2924\"],\n  [\"99\u00d737=3663\", \"85\u00d795=8075\"],\n  [\"81\u00d753=4293\", \"13\u00d788=1144\"],\n  [\"88\u00d795=8360\", \"79\u00d750=3950\"],\n  [\"55\u00d782=4510\", \"97\u00d730=2910\"],\n  [\"18\u00d725=450\", \"40\u00d767=2680\"],\n  [\"51\u00d777=3927\", \"32\u00d795=3040\"],\n  [\"13\u00d732=416\", \"59\u00d766=3894\"],\n  [\"48\u00d764=3072\", \"50\u00d711=550\"],\n  [\"78\u00d711=858\", \"97\u00d794=9118\"],\n  [\"64\u00d746=2944\", \"81\u00d769=5589\"],\n  [\"90\u00d754=4860\", \"52\u00d743=2236\"],\n  [\"75\u00d762=4650\", \"16\u00d776=1216\"],\n  [\"57\u00d774=4218\", \"21\u00d758=1218\"],\n  [\"73\u00d763=4599\", \"94\u00d759=5546\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and the 25 \"two-digit x two-digit\" answer cells\n# with the updated values from the new day's worksheet. Every old value\n# in this document is unique, so a simple exact, case-sensitive\n# Find/Replace (replace first occurrence) for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-11 Wednesday\", \"2024-12-12 Thursday\"),\n    @(\"44\u00d728=1232\", \"74\u00d740=2960\"),\n    @(\"78\u00d728=2184\", \"23\u00d724=552\"),\n    @(\"22\u00d795=2090\", \"87\u00d718=1566\"),\n    @(\"41\u00d734=1394\", \"73\u00d775=5475\"),\n    @(\"19\u00d723=437\", \"83\u00d778=6474\"),\n    @(\"66\u00d723=1518\", \"94\u00d717=1598\"),\n    @(\"39\u00d726=1014\", \"40\u00d779=3160\"),\n    @(\"33\u00d775=2475\", \"13\u00d791=1183\"),\n    @(\"98\u00d756=5488\", \"86\u00d760=5160\"),\n    @(\"79\u00d721=1659\", \"39\u00d728=1092\"),\n    @(\"48\u00d769=3312\", \"34\u00d786=2924\"),\n    @(\"99\u00d737=3663\", \"85\u00d795=8075\"),\n    @(\"81\u00d753=4293\", \"13\u00d788=1144\"),\n    @(\"88\u00d795=8360\", \"79\u00d750=3950\"),\n    @(\"55\u00d782=4510\", \"97\u00d730=2910\"),\n    @(\"18\u00d725=450\", \"40\u00d767=2680\"),\n    @(\"51\u00d777=3927\", \"32\u00d795=3040\"),\n    @(\"13\u00d732=416\", \"59\u00d766=3894\"),\n    @(\"48\u00d764=3072\", \"50\u00d711=550\"),\n    @(\"78\u00d711=858\", \"97\u00d794=9118\"),\n    @(\"64\u00d746=2944\", \"81\u00d769=5589\"),\n    @(\"90\u00d754=4860\", \"52\u00d743=2236\"),\n    @(\"75\u00d762=4650\", \"16\u00d776=1216\"),\n    @(\"57\u00d774=4218\", \"21\u00d758=1218\"),\n    @(\"73\u00d763=4599\", \"94\u00d759=5546\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$newText, 2)\n}\n"}
